$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# Update the "Days" legend label to reflect the actual start date
$ws.Range("H3").Value = "Days. Started on 7.8.23"

# Row 5 - 1. Initial Planning
$ws.Range("D5").Value = 4
$ws.Range("F5").Value = 3

# Row 6 - 1.1 Early Planning
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Row 7 - 1.2 Work Division
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1

# Row 8 - 1.3 Initial Research
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 4

# Row 9 - 2. Project Plan Doc
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 1

# Row 10 - 2.1 Introduction
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 1

# Row 11 - 2.2 Activity Defintion
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 4

# Row 12 - 2.3 Gantt Chart
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 15
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 14
$ws.Range("G12").Value = 1

# Row 13 - 3. Software Design Doc
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 5
$ws.Range("G13").Value = 1

# Row 14 - 3.1 Software Vision Doc
$ws.Range("C14").Value = 12
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 1

# Row 15 - 3.2 System Requirements
$ws.Range("C15").Value = 12
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 1

# Row 16 - 3.3 Components and Software
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 1

# Row 17 - 3.4 User Interface
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 1

# Update selected/active cell to reflect where the author finished editing
$ws.Range("Z17").Select()
